# sprint-backlog update for sprint 5/6:
# record the actual execution status for tasks T29 (J), T30 (K) and T32 (C)
# on the "execution" sheet, then leave the "execution" sheet active/selected
# (it was "burndown" before).

$wb = $excel.ActiveWorkbook

$execution = $wb.Worksheets.Item("execution")

# Task T29 ("J") — day 4 now shows 2 points remaining, day 5 shows 1.
$execution.Range("H4").Value = "J:2"
$execution.Range("I4").Value = "J:1"

# Task T30 ("K") — day 5 shows 6 points remaining.
$execution.Range("I5").Value = "K:6"

# Task T32 ("C") — day 3 shows 4 points remaining, day 5 shows 2.
$execution.Range("G7").Value = "C:4"
$execution.Range("I7").Value = "C:2"

# Make "execution" the active sheet/tab, with H4 selected.
$execution.Activate()
$execution.Range("H4").Select()
